# Roche Phase 3 NCT04320615.xlsx -- "First cut for endpoint and objectives"
# Adds a new worksheet "studyDesignOE" (Objectives & Endpoints) at the end of
# the workbook, populates its header + data rows, and makes it the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new worksheet after the last existing sheet (studyDesignPopulations)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "studyDesignOE"

# ---------------------------------------------------------------------------
# 2. Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 43.166666666666664
$ws.Columns.Item(2).ColumnWidth = 17.998697916666668
$ws.Columns.Item(3).ColumnWidth = 58.666666666666664
$ws.Columns.Item(4).ColumnWidth = 30.666666666666668
$ws.Columns.Item(5).ColumnWidth = 30.666666666666668

# ---------------------------------------------------------------------------
# 3. Header row
# ---------------------------------------------------------------------------
$ws.Range("A1").Value2 = "objectiveDescription"
$ws.Range("B1").Value2 = "objectiveLevel"
$ws.Range("C1").Value2 = "endpointDescription"
$ws.Range("D1").Value2 = "endpointPurposeDescription"
$ws.Range("E1").Value2 = "endpointLevel"

# ---------------------------------------------------------------------------
# 4. Data rows
# ---------------------------------------------------------------------------
$ws.Range("A2").Value2 = "The primary efficacy objective for this study is to evaluate the efficacy of TCZ compared with placebo in combination with SOC for the treatment of severe COVID-19 pneumonia"
$ws.Range("B2").Value2 = "Study Primary Objective"
$ws.Range("C2").Value2 = "Clinical status assessed using a 7-category ordinal scale at Day 28"
$ws.Range("E2").Value2 = "Primary Endpoint"

$ws.Range("A3").Value2 = "The secondary efficacy objective for this study is to evaluate the efficacy of TCZ compared with placebo in combination with SOC for the treatment of severe COVID-19 pneumonia"
$ws.Range("B3").Value2 = "Study Secondary Objective"
$ws.Range("C3").Value2 = "Time to clinical improvement (TTCI) defined as a National Early Warning Score 2 (NEWS2) of <=2 maintained for 24 hours"
$ws.Range("E3").Value2 = "Secondary Enpoint"

$ws.Range("C4").Value2 = "Time to improvement of at least 2 categories relative to baseline on a 7-category ordinal scale of clinical status"
$ws.Range("E4").Value2 = "Secondary Enpoint"

$ws.Range("C5").Value2 = "Incidence of mechanical ventilation"
$ws.Range("E5").Value2 = "Secondary Enpoint"

$ws.Range("C6").Value2 = "Ventilator-free days to Day 28"
$ws.Range("E6").Value2 = "Secondary Enpoint"

$ws.Range("C7").Value2 = "Incidence of intensive care unit (ICU) stay"
$ws.Range("E7").Value2 = "Secondary Enpoint"

$ws.Range("C8").Value2 = "Duration of ICU stay"
$ws.Range("E8").Value2 = "Secondary Enpoint"

$ws.Range("C9").Value2 = "Time to clinical failure, defined as the time to death, mechanical ventilation, ICU admission, or withdrawal (whichever occurs first). For patients entering the study already in ICU or on mechanical ventilation, clinical failure is defined as a one-category worsening on the ordinal scale, withdrawal or death."
$ws.Range("E9").Value2 = "Secondary Enpoint"

$ws.Range("C10").Value2 = "Mortality rate at Days 7, 14, 21, 28, and 60"
$ws.Range("E10").Value2 = "Secondary Enpoint"

$ws.Range("C11").Value2 = "Time to hospital discharge or “ready for discharge” (as evidenced by normal body temperature and respiratory rate, and stable oxygen saturation on ambient air or <= 2L supplemental oxygen)"
$ws.Range("E11").Value2 = "Secondary Enpoint"

$ws.Range("C12").Value2 = "Time to recovery, defined as discharged or “ready for discharge” (as evidenced by normal body temperature and respiratory rate, and stable oxygen saturation on ambient air or <= 2L supplemental oxygen); OR, in a non-ICU hospital ward (or “ready for hospital ward”) not requiring supplemental oxygen"
$ws.Range("E12").Value2 = "Secondary Enpoint"

$ws.Range("C13").Value2 = "Duration of supplemental oxygen"
$ws.Range("E13").Value2 = "Secondary Enpoint"

# ---------------------------------------------------------------------------
# 5. Formatting: header (bold, left/top aligned, wrapped) and body
#    (left/top aligned, wrapped) across the full used range A1:P36
# ---------------------------------------------------------------------------
$hdr = $ws.Range("A1:E1")
$hdr.HorizontalAlignment = -4131
$hdr.VerticalAlignment = -4160
$hdr.WrapText = $true
$hdr.Font.Bold = $true

$body = $ws.Range("A2:P36")
$body.HorizontalAlignment = -4131
$body.VerticalAlignment = -4160
$body.WrapText = $true

$rowF = $ws.Range("F1:P1")
$rowF.HorizontalAlignment = -4131
$rowF.VerticalAlignment = -4160
$rowF.WrapText = $true

# ---------------------------------------------------------------------------
# 6. Row heights (matching Excel's computed wrap heights for the content)
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 17
$ws.Rows.Item(2).RowHeight = 68
$ws.Rows.Item(3).RowHeight = 68
$ws.Rows.Item(4).RowHeight = 34
$ws.Rows.Item(5).RowHeight = 17
$ws.Rows.Item(6).RowHeight = 17
$ws.Rows.Item(7).RowHeight = 17
$ws.Rows.Item(8).RowHeight = 17
$ws.Rows.Item(9).RowHeight = 85
$ws.Rows.Item(10).RowHeight = 17
$ws.Rows.Item(11).RowHeight = 51
$ws.Rows.Item(12).RowHeight = 85
$ws.Rows.Item(13).RowHeight = 17

# ---------------------------------------------------------------------------
# 7. View: zoom, selection, and make this the active/selected tab
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A1:E1048576").Select()
$excel.ActiveWindow.Zoom = 130
